$d = $word.ActiveDocument

# ---------------------------------------------------------------------
# 1. Remove the "Meta description: ..." paragraph that currently
#    follows the Heading1 title paragraph.
# ---------------------------------------------------------------------
$metaPara = $d.Paragraphs.Item(2)
$metaPara.Range.Delete()

# ---------------------------------------------------------------------
# 2. Replace the final "Prompt: ..." paragraph with two paragraphs:
#      a) a new bold paragraph repeating the page title, and
#      b) the same paragraph (now italic) but holding the new
#         meta-description copy instead of the old image prompt.
#    InsertXML on a collapsed range replaces the whole paragraph that
#    contains the insertion point, so we supply both paragraphs' XML
#    together to end up with one extra paragraph overall.
# ---------------------------------------------------------------------
$w = "xmlns:w='http://schemas.openxmlformats.org/wordprocessingml/2006/main'"

$titlePara = "<w:p $w><w:r/><w:r><w:rPr><w:b/></w:rPr><w:t>Play Eye of Horus Megaways Free: A Modernized Ancient Egypt Slot</w:t></w:r></w:p>"
$descPara  = "<w:p $w><w:r/><w:r><w:rPr><w:i/></w:rPr><w:t>Read our review of Eye of Horus Megaways and play it for free. Experience the modern version of an Ancient Egypt slot with up to 15,625 ways to win.</w:t></w:r></w:p>"

$lastPara = $d.Paragraphs.Item($d.Paragraphs.Count)
$insertRange = $lastPara.Range
$insertRange.Collapse(1)
$insertRange.InsertXML($titlePara + $descPara)
